# Generate Report for Handback
# Update "Latest Handback DateTime" (column K) for the c1298a7c-... row (row 2)
# on both the "zh-cn" and "de-de" sheets, reflecting freshly generated handback
# timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("K2").Value = "2016-11-02 04:58:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-11-02 04:59:12"
